# Auto-generated Excel COM-interop script
# Applies the "Updated cryptos list" data refresh (Fri Sep 22 23:44:59 UTC 2023)
# to the cryptocurrency listing sheet: refreshed prices / 1h volume percentages,
# plus a couple of rank swaps and one coin being replaced in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain-looking numbers (e.g. "211.76", "1.00") as
# literal text, matching the site's display strings exactly (including
# trailing zeros). Force those specific cells to Text format first so Excel's
# automatic type inference does not silently convert them to numeric values
# and strip formatting (e.g. "1.00" -> 1, "2.20" -> 2.2).
$textCells = @(
    "D5",
    "D6",
    "D10",
    "D11",
    "D16",
    "D19",
    "D25",
    "D29",
    "D33",
    "D41",
    "D42",
    "D43",
    "D46",
    "D47",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.660.30"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.596.85"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "211.76"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "0.515"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "19.53"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "0.0837"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "1.821.15"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "1.624.59"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D16").Value = "64.45"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "26.633.19"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "0.0₃0733"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "208.96"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "145.30"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("D29").Value = "15.27"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "0.659"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("D35").Value = "1.277.48"
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "5.48"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "64.59"
$ws.Range("E42").Value = "  +3.28%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.20"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.733.73"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "0.915"
$ws.Range("E46").Value = "  +9.24%  "
$ws.Range("D47").Value = "90.03"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  +4.68%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.10%  "
